# Train the model for SVM
# Insert a "Category" column (new column C), shifting the old "Body" column to D,
# and populate it with the classification produced for each email.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C ("Body"), shifting it to D.
$ws.Columns.Item(3).Insert()

# Header for the new column, matching the style used by the other header cells.
$ws.Cells.Item(1, 3).Value = "Category"
$ws.Cells.Item(1, 3).Style = $ws.Cells.Item(1, 1).Style

# Column widths: new Category column, and the (shifted) former Body column.
$ws.Columns.Item(3).ColumnWidth = 26
$ws.Columns.Item(4).ColumnWidth = 9.140625

$categories = @(
    "Health",
    "Automobile",
    "Automobile",
    "Building",
    "Building",
    "Health",
    "Health",
    "Automobile",
    "Building",
    "Automobile",
    "Automobile",
    "Health",
    "Health",
    "Automobile",
    "Building",
    "Health",
    "Building",
    "Automobile",
    "Health",
    "Health",
    "Building",
    "Automobile",
    "Building",
    "Automobile",
    "Health",
    "Automobile",
    "Health",
    "Health",
    "Health",
    "Building",
    "Automobile",
    "Automobile",
    "Building",
    "Automobile",
    "Health",
    "Health",
    "Health",
    "Building",
    "Health",
    "Automobile",
    "Health",
    "Automobile",
    "Building",
    "Building",
    "Building",
    "Health",
    "Automobile",
    "Health",
    "Building",
    "Health",
    "Building",
    "Automobile",
    "Health",
    "Automobile",
    "Building",
    "Health",
    "Health",
    "Building",
    "Health",
    "Building"
)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $categories[$i]
}

# Scroll/selection state, matching where the editing session ended up.
$ws.Application.ActiveWindow.ScrollRow = 45
$ws.Range("C61").Select()
